$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date serial that was updated from
# 2023-10-05 (45204) to 2023-10-08 (45207) for every data row (rows 2-89).
$ws.Range("C2:C89").Value = 45207
